$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Sending cluster" label in row 2 from "ECs" to "MuSCs"
$ws.Range("A2").Value = "MuSCs"

# Update recomputed TPM-derived values
$ws.Range("G2").Value = 0.2195956666666667
$ws.Range("H2").Value = 0.658787
$ws.Range("Q2").Value = 2.529134092797555
$ws.Range("R2").Value = 22.762206835178
